$wb = $excel.ActiveWorkbook

# ALC row 69
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(69, 8).Value = 1604.1666
$ws.Cells.Item(69, 10).Value = 1540.909
$ws.Cells.Item(69, 12).Value = 4622.727000000001
$ws.Cells.Item(69, 14).Value = -6370.727000000001

# ALC row 72
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(72, 8).Value = 1604.1666
$ws.Cells.Item(72, 10).Value = 1540.909
$ws.Cells.Item(72, 12).Value = 13868.181
$ws.Cells.Item(72, 14).Value = -22604.181

# ALC row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 5940.5713
$ws.Cells.Item(86, 9).Value = 1264.6
$ws.Cells.Item(86, 11).Value = 1264.6
$ws.Cells.Item(86, 13).Value = -141.5999999999999

# ALC row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(89, 8).Value = 5940.5713
$ws.Cells.Item(89, 9).Value = 1264.6
$ws.Cells.Item(89, 11).Value = 6323
$ws.Cells.Item(89, 13).Value = -707

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 2713.6924
$ws.Cells.Item(138, 9).Value = 2217.9092
$ws.Cells.Item(138, 10).Value = 2846.7073
$ws.Cells.Item(138, 11).Value = 6653.7276
$ws.Cells.Item(138, 12).Value = 8540.1219
$ws.Cells.Item(138, 13).Value = -1513.7276
$ws.Cells.Item(138, 14).Value = -18820.1219

# ARM row 5
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 138
$ws.Cells.Item(5, 9).Value = 175
$ws.Cells.Item(5, 10).Value = 101
$ws.Cells.Item(5, 11).Value = 175
$ws.Cells.Item(5, 12).Value = 101
$ws.Cells.Item(5, 13).Value = -63
$ws.Cells.Item(5, 14).Value = -325

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6850.159
$ws.Cells.Item(32, 9).Value = 5186.65
$ws.Cells.Item(32, 10).Value = 17940.223
$ws.Cells.Item(32, 11).Value = 5186.65
$ws.Cells.Item(32, 12).Value = 17940.223
$ws.Cells.Item(32, 13).Value = -4899.65
$ws.Cells.Item(32, 14).Value = -18514.223

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 2085114.6
$ws.Cells.Item(63, 9).Value = 2055.9092
$ws.Cells.Item(63, 10).Value = 7813526
$ws.Cells.Item(63, 11).Value = 2055.9092
$ws.Cells.Item(63, 12).Value = 7813526
$ws.Cells.Item(63, 13).Value = -1369.9092
$ws.Cells.Item(63, 14).Value = -7814898

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(66, 8).Value = 2085114.6
$ws.Cells.Item(66, 9).Value = 2055.9092
$ws.Cells.Item(66, 10).Value = 7813526
$ws.Cells.Item(66, 11).Value = 10279.546
$ws.Cells.Item(66, 12).Value = 39067630
$ws.Cells.Item(66, 13).Value = -6847.546
$ws.Cells.Item(66, 14).Value = -39074494

# BSM row 4
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 138
$ws.Cells.Item(4, 9).Value = 175
$ws.Cells.Item(4, 10).Value = 101
$ws.Cells.Item(4, 11).Value = 175
$ws.Cells.Item(4, 12).Value = 101
$ws.Cells.Item(4, 13).Value = -60
$ws.Cells.Item(4, 14).Value = -331

# BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1267.18
$ws.Cells.Item(20, 9).Value = 1175.7812
$ws.Cells.Item(20, 10).Value = 1429.6666
$ws.Cells.Item(20, 11).Value = 1175.7812
$ws.Cells.Item(20, 12).Value = 1429.6666
$ws.Cells.Item(20, 13).Value = -928.7811999999999
$ws.Cells.Item(20, 14).Value = -1923.6666

# BSM row 35
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(35, 8).Value = 24628.8
$ws.Cells.Item(35, 10).Value = 24628.8
$ws.Cells.Item(35, 12).Value = 24628.8
$ws.Cells.Item(35, 14).Value = -25248.8

# CRP row 19
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(19, 8).Value = 277
$ws.Cells.Item(19, 9).Value = 246.25
$ws.Cells.Item(19, 10).Value = 400
$ws.Cells.Item(19, 11).Value = 246.25
$ws.Cells.Item(19, 12).Value = 400
$ws.Cells.Item(19, 13).Value = -76.25
$ws.Cells.Item(19, 14).Value = -740

# CRP row 24
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(24, 8).Value = 277
$ws.Cells.Item(24, 9).Value = 246.25
$ws.Cells.Item(24, 10).Value = 400
$ws.Cells.Item(24, 11).Value = 246.25
$ws.Cells.Item(24, 12).Value = 400
$ws.Cells.Item(24, 13).Value = -76.25
$ws.Cells.Item(24, 14).Value = -740

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3440.186
$ws.Cells.Item(31, 9).Value = 2509.6667
$ws.Cells.Item(31, 10).Value = 3938.6785
$ws.Cells.Item(31, 11).Value = 2509.6667
$ws.Cells.Item(31, 12).Value = 3938.6785
$ws.Cells.Item(31, 13).Value = -2214.6667
$ws.Cells.Item(31, 14).Value = -4528.6785

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 3440.186
$ws.Cells.Item(34, 9).Value = 2509.6667
$ws.Cells.Item(34, 10).Value = 3938.6785
$ws.Cells.Item(34, 11).Value = 2509.6667
$ws.Cells.Item(34, 12).Value = 3938.6785
$ws.Cells.Item(34, 13).Value = -2307.6667
$ws.Cells.Item(34, 14).Value = -4342.6785

# CRP row 41
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(41, 8).Value = 19971.428
$ws.Cells.Item(41, 10).Value = 19971.428
$ws.Cells.Item(41, 12).Value = 19971.428
$ws.Cells.Item(41, 14).Value = -20827.428

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 958.7593000000001
$ws.Cells.Item(5, 9).Value = 559.0540999999999
$ws.Cells.Item(5, 11).Value = 1677.1623
$ws.Cells.Item(5, 13).Value = -1565.1623

# CUL row 86
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(86, 8).Value = 2611
$ws.Cells.Item(86, 9).Value = 2222
$ws.Cells.Item(86, 10).Value = 3000
$ws.Cells.Item(86, 11).Value = 6666
$ws.Cells.Item(86, 12).Value = 9000
$ws.Cells.Item(86, 13).Value = -5480
$ws.Cells.Item(86, 14).Value = -11372

# CUL row 89
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(89, 8).Value = 2611
$ws.Cells.Item(89, 9).Value = 2222
$ws.Cells.Item(89, 10).Value = 3000
$ws.Cells.Item(89, 11).Value = 19998
$ws.Cells.Item(89, 12).Value = 27000
$ws.Cells.Item(89, 13).Value = -14070
$ws.Cells.Item(89, 14).Value = -38856

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(135, 8).Value = 958.7593000000001
$ws.Cells.Item(135, 9).Value = 559.0540999999999
$ws.Cells.Item(135, 11).Value = 5031.4869
$ws.Cells.Item(135, 13).Value = -2496.4869

# GSM row 46
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(46, 8).Value = 24720
$ws.Cells.Item(46, 10).Value = 26293.334
$ws.Cells.Item(46, 12).Value = 26293.334
$ws.Cells.Item(46, 14).Value = -26605.334

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 3474.2
$ws.Cells.Item(80, 9).Value = 3335.5
$ws.Cells.Item(80, 10).Value = 3566.6667
$ws.Cells.Item(80, 11).Value = 3335.5
$ws.Cells.Item(80, 12).Value = 3566.6667
$ws.Cells.Item(80, 13).Value = -2337.5
$ws.Cells.Item(80, 14).Value = -5562.6667

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(83, 8).Value = 3474.2
$ws.Cells.Item(83, 9).Value = 3335.5
$ws.Cells.Item(83, 10).Value = 3566.6667
$ws.Cells.Item(83, 11).Value = 16677.5
$ws.Cells.Item(83, 12).Value = 17833.3335
$ws.Cells.Item(83, 13).Value = -11685.5
$ws.Cells.Item(83, 14).Value = -27817.3335

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 3534.5
$ws.Cells.Item(22, 10).Value = 1433.3334
$ws.Cells.Item(22, 12).Value = 1433.3334
$ws.Cells.Item(22, 14).Value = -2023.3334

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(27, 8).Value = 3534.5
$ws.Cells.Item(27, 10).Value = 1433.3334
$ws.Cells.Item(27, 12).Value = 1433.3334
$ws.Cells.Item(27, 14).Value = -1647.3334

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 418089.28
$ws.Cells.Item(132, 9).Value = 525460.5
$ws.Cells.Item(132, 10).Value = 6499.6665
$ws.Cells.Item(132, 11).Value = 1576381.5
$ws.Cells.Item(132, 12).Value = 19498.9995
$ws.Cells.Item(132, 13).Value = -1573851.5
$ws.Cells.Item(132, 14).Value = -24558.9995

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1092.5769
$ws.Cells.Item(132, 9).Value = 557.3333
$ws.Cells.Item(132, 10).Value = 3340.6
$ws.Cells.Item(132, 11).Value = 1671.9999
$ws.Cells.Item(132, 12).Value = 10021.8
$ws.Cells.Item(132, 13).Value = 858.0001
$ws.Cells.Item(132, 14).Value = -15081.8
